$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3: swap stimulus from img_5jy9c.png to img_8dmpq.png, with the values
# that previously belonged to row 28 (the img_8dmpq.png record).
$ws.Range("L3").Value = "stimuli/img_8dmpq.png"
$ws.Range("M3").Value = 30.65909090909091
$ws.Range("N3").Value = 24.11363636363636
$ws.Range("O3").Value = 27.38636363636364
$ws.Range("P3").Value = 44
$ws.Range("Q3").Value = 2
$ws.Range("R3").Value = 2
$ws.Range("S3").Value = 2
$ws.Range("T3").Value = 2
$ws.Range("U3").Value = 2
$ws.Range("V3").Value = 2

# Row 16: rename catch stimulus file
$ws.Range("L16").Value = "stimuli/catch_01.jpg"

# Row 28: swap stimulus from img_8dmpq.png to img_5jy9c.png, with the values
# that previously belonged to row 3 (the img_5jy9c.png record).
$ws.Range("L28").Value = "stimuli/img_5jy9c.png"
$ws.Range("M28").Value = 87.37209302325581
$ws.Range("N28").Value = 79.18604651162791
$ws.Range("O28").Value = 83.27906976744185
$ws.Range("P28").Value = 43
$ws.Range("Q28").Value = 10
$ws.Range("R28").Value = 10
$ws.Range("S28").Value = 10
$ws.Range("T28").Value = 10
$ws.Range("U28").Value = 9
$ws.Range("V28").Value = 10
